# Applies the "Created and validated the reference sequence set" edit:
#   - Column B ("Class") labels are renamed from the old short codes to the
#     new IFNL-prefixed labels:
#       AMPH      -> IFNL-Amphibian
#       REPTILE1  -> IFNLa-Reptile
#       REPTILE2  -> IFNL-Reptile
#       AVES      -> IFNL-Aves
#       MAMM (Clade B) -> IFNL-Mammal
#       MAMM (Clade A) -> IFNLa-Mammal
#   - Column D ("GeneID") rows 10-13 had a placeholder "MAMM" value that is
#     corrected to "MAMM1" (matching the rest of that GeneID group).
#   - F4 (previously empty "Species" cell for a REPTILE2 row) is filled in
#     with "Pelodiscus sinensis".
#   - The sheet/window view state (scroll position + selection) is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("A+B")

# --- Column B ("Class") relabeling -----------------------------------
$ws.Range("B2").Value        = "IFNL-Amphibian"
$ws.Range("B3").Value        = "IFNLa-Reptile"
$ws.Range("B4:B6").Value     = "IFNL-Reptile"
$ws.Range("B7:B9").Value     = "IFNL-Aves"
$ws.Range("B10:B66").Value   = "IFNL-Mammal"
$ws.Range("B67:B108").Value  = "IFNLa-Mammal"
$ws.Range("B109:B110").Value = "IFNL-Mammal"

# --- Column D ("GeneID") data correction ------------------------------
$ws.Range("D10:D13").Value = "MAMM1"

# --- New species value for previously-empty cell ----------------------
$ws.Range("F4").Value = "Pelodiscus sinensis"

# --- View state (scroll position / active selection) ------------------
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("E36").Select()

